$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-file names added to the "Piotr Bistyga" block (columns K:M)
# alongside matching dates (K) and line counts (M).
$rows = @(
    @{ Row = 22; Date = 45793; File = "WykonaneBadaniaControllerTests.cs"; Lines = 44 },
    @{ Row = 23; Date = 45793; File = "PacjentControllerTests.cs";         Lines = 41 },
    @{ Row = 24; Date = 45793; File = "RecepcjonistkaControllerTests.cs";  Lines = 45 },
    @{ Row = 25; Date = 45793; File = "WizytaControllerTests.cs";          Lines = 60 }
)

foreach ($r in $rows) {
    # Copy the existing date-column formatting (column K already uses it
    # two rows up) instead of re-creating a numeric format, so no new
    # style entries are introduced.
    $ws.Cells.Item(20, 11).Copy()
    $ws.Cells.Item($r.Row, 11).PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($r.Row, 11).Value = $r.Date
    $ws.Cells.Item($r.Row, 12).Value = $r.File
    $ws.Cells.Item($r.Row, 13).Value = $r.Lines
}

# Restore the selection to where the editor last left it.
$ws.Range("M26").Select()
